$wb = $excel.ActiveWorkbook

# Rename the "Add_Round" sheet to "Add_Round_Player1"
$roundSheet = $wb.Worksheets.Item("Add_Round")
$roundSheet.Name = "Add_Round_Player1"

# Update the selection on "Simple_Out_Flight_Player1" (G11 -> G8) without
# leaving it as the active tab.
$simpleSheet = $wb.Worksheets.Item("Simple_Out_Flight_Player1")
$simpleSheet.Range("G8").Select()

# Finally activate "Add_Round_Player1" and move its selection (E10 -> F9);
# this becomes the workbook's active/selected tab.
$roundSheet.Activate()
$roundSheet.Range("F9").Select()
